$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'23.342.88"
$ws.Range("E2").Value = "  -0.54%  "
$ws.Range("D3").Value = "'1.625.22"
$ws.Range("E3").Value = "  -0.87%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").Value = "'1.003"
$ws.Range("D6").Value = "'303.38"
$ws.Range("E6").Value = "  -0.74%  "
$ws.Range("D7").Value = "'0.3747"
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "'0.3622"
$ws.Range("E8").Value = "  -0.38%  "
$ws.Range("D9").Value = "'51.38"
$ws.Range("E9").Value = "  -1.22%  "
$ws.Range("D10").Value = "'0.08150"
$ws.Range("E10").Value = "  +0.09%  "
$ws.Range("D11").Value = "'1.224"
$ws.Range("E11").Value = "  -2.90%  "
$ws.Range("D12").Value = "'1.002"
$ws.Range("E12").Value = "  +0.19%  "
$ws.Range("D13").Value = "'22.23"
$ws.Range("E13").Value = "  -3.15%  "
$ws.Range("D14").Value = "'6.473"
$ws.Range("E14").Value = "  -2.31%  "
$ws.Range("D15").Value = "'0.00001237"
$ws.Range("E15").Value = "  -2.99%  "
$ws.Range("D16").Value = "'7.285"
$ws.Range("E16").Value = "  -0.83%  "
$ws.Range("D17").Value = "'1.626.34"
$ws.Range("E17").Value = "  -0.28%  "
$ws.Range("D18").Value = "'93.72"
$ws.Range("E18").Value = "  -0.80%  "
$ws.Range("D19").Value = "'0.06953"
$ws.Range("E19").Value = "  +0.66%  "
$ws.Range("D20").Value = "'17.46"
$ws.Range("E20").Value = "  -4.05%  "
$ws.Range("D21").Value = "'6.507"
$ws.Range("E21").Value = "  -0.44%  "
$ws.Range("D22").Value = "'1.003"
$ws.Range("E22").Value = "  +0.28%  "
$ws.Range("D23").Value = "'12.53"
$ws.Range("E23").Value = "  -1.83%  "
$ws.Range("D24").Value = "'23.334.51"
$ws.Range("E24").Value = "  -0.56%  "
$ws.Range("D25").Value = "'3.131"
$ws.Range("E25").Value = "  +1.00%  "
$ws.Range("D26").Value = "'2.447"
$ws.Range("E26").Value = "  +1.18%  "
$ws.Range("D27").Value = "'21.22"
$ws.Range("E27").Value = "  -0.02%  "
$ws.Range("D28").Value = "'150.66"
$ws.Range("E28").Value = "  -0.14%  "
$ws.Range("D29").Value = "'5.267"
$ws.Range("E29").Value = "  -1.21%  "
$ws.Range("D30").Value = "'132.67"
$ws.Range("E30").Value = "  -2.76%  "
$ws.Range("D31").Value = "'1.798.08"
$ws.Range("E31").Value = "  -0.81%  "
$ws.Range("D32").Value = "'2.226"
$ws.Range("E32").Value = "  -3.36%  "
$ws.Range("D33").Value = "'6.749"
$ws.Range("E33").Value = "  -0.12%  "
$ws.Range("D34").Value = "'1.027"
$ws.Range("E34").Value = "  +6.98%  "
$ws.Range("D35").Value = "'10.65"
$ws.Range("E35").Value = "  +3.06%  "
$ws.Range("D36").Value = "'0.02749"
$ws.Range("E36").Value = "  -3.17%  "
$ws.Range("D37").Value = "'0.2491"
$ws.Range("E37").Value = "  -1.41%  "
$ws.Range("D38").Value = "'0.08765"
$ws.Range("E38").Value = "  -0.47%  "
$ws.Range("D39").Value = "'0.07105"
$ws.Range("E39").Value = "  -2.67%  "
$ws.Range("D40").Value = "'5.957"
$ws.Range("E40").Value = "  -2.63%  "
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").Value = "'0.6971"
$ws.Range("E41").Value = "  -1.64%  "
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").Value = "'1.333"
$ws.Range("E42").Value = "  -3.15%  "
$ws.Range("D43").Value = "'15.88"
$ws.Range("E43").Value = "  -1.69%  "
$ws.Range("D44").Value = "'12.05"
$ws.Range("E44").Value = "  -3.27%  "
$ws.Range("D45").Value = "'0.6483"
$ws.Range("E45").Value = "  -0.97%  "
$ws.Range("D46").Value = "'1.002"
$ws.Range("E46").Value = "  +0.17%  "
$ws.Range("D47").Value = "'2.270"
$ws.Range("E47").Value = "  -2.64%  "
$ws.Range("E48").Value = "  -1.18%  "
$ws.Range("D49").Value = "'0.07961"
$ws.Range("E49").Value = "  -0.08%  "
$ws.Range("D50").Value = "'1.188"
$ws.Range("E50").Value = "  -1.51%  "
$ws.Range("D51").Value = "'125.16"
$ws.Range("E51").Value = "  -2.81%  "
